# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.034748368925986)
    3 = @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.21235251628842)
    4 = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 1, 4.266567469086664)
    5 = @(1.505614041169197, 1.65323645889881, 16.98373111632243, 6.48142807727062, 1, 26.62400969366105)
    6 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
